# Apply cryptos list update per diff (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.337.83'
$ws.Range('E2').Value = '  -2.87%  '
$ws.Range('D3').Value = '3.498.61'
$ws.Range('E3').Value = '  -2.37%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '''608.27'
$ws.Range('E5').Value = '  +4.05%  '
$ws.Range('D6').Value = '''185.56'
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').Value = '''0.626'
$ws.Range('E7').Value = '  +0.40%  '
$ws.Range('D8').Value = '''0.999'
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.210'
$ws.Range('E9').Value = '  -3.35%  '
$ws.Range('D10').Value = '''0.653'
$ws.Range('E10').Value = '  -0.57%  '
$ws.Range('D11').Value = '''53.26'
$ws.Range('E11').Value = '  -2.79%  '
$ws.Range('D12').Value = '''0.0000305'
$ws.Range('E12').Value = '  -4.42%  '
$ws.Range('D13').Value = '''9.63'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').Value = '4.061.67'
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').Value = '''613.02'
$ws.Range('E15').Value = '  +8.20%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = '''12.71'
$ws.Range('E16').Value = '  +2.33%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').Value = '''18.98'
$ws.Range('E17').Value = '  -1.80%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '69.388.16'
$ws.Range('E18').Value = '  -2.63%  '
$ws.Range('D19').Value = '3.497.53'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('D21').Value = '''0.988'
$ws.Range('E21').Value = '  -2.62%  '
$ws.Range('D22').Value = '''17.57'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = '''105.07'
$ws.Range('E23').Value = '  +10.69%  '
$ws.Range('D24').Value = '''4.65'
$ws.Range('E24').Value = '  +1.45%  '
$ws.Range('D25').Value = '''5.02'
$ws.Range('E25').Value = '  -1.10%  '
$ws.Range('D26').Value = '''3.01'
$ws.Range('E26').Value = '  +1.57%  '
$ws.Range('D27').Value = '''10.91'
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('D28').Value = '''9.90'
$ws.Range('E28').Value = '  +7.65%  '
$ws.Range('D29').Value = '''33.76'
$ws.Range('E29').Value = '  +3.24%  '
$ws.Range('D30').Value = '''6.99'
$ws.Range('E30').Value = '  -4.58%  '
$ws.Range('D31').Value = '''12.45'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('D32').Value = '''0.116'
$ws.Range('E32').Value = '  -0.34%  '
$ws.Range('D33').Value = '''63.41'
$ws.Range('E33').Value = '  -1.33%  '
$ws.Range('D34').Value = '''3.72'
$ws.Range('E34').Value = '  +13.64%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '''3.13'
$ws.Range('E35').Value = '  -8.59%  '
$ws.Range('B36').Value = 'Dai'
$ws.Range('C36').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '''525.74'
$ws.Range('E37').Value = '  -5.26%  '
$ws.Range('D38').Value = '''0.396'
$ws.Range('E38').Value = '  -6.20%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.546.15'
$ws.Range('E39').Value = '  +0.92%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''3.57'
$ws.Range('E40').Value = '  +3.34%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').Value = '''36.56'
$ws.Range('E41').Value = '  -3.41%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '''0.140'
$ws.Range('E42').Value = '  +2.73%  '
$ws.Range('D43').Value = '0.0₃0765'
$ws.Range('E43').Value = '  -5.48%  '
$ws.Range('D44').Value = '''0.0459'
$ws.Range('E44').Value = '  +1.96%  '
$ws.Range('D45').Value = '''2.94'
$ws.Range('E45').Value = '  -0.65%  '
$ws.Range('D46').Value = '''0.143'
$ws.Range('E46').Value = '  +3.85%  '
$ws.Range('D47').Value = '''3.33'
$ws.Range('E47').Value = '  -3.89%  '
$ws.Range('D48').Value = '''8.88'
$ws.Range('E48').Value = '  -5.62%  '
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').Value = '''132.03'
$ws.Range('E50').Value = '  -2.71%  '
$ws.Range('D51').Value = '''1.35'
$ws.Range('E51').Value = '  -9.01%  '
